$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new columns before column D -----------------------------
# This shifts the existing headers (old D:Q, "N_STO = 7" .. "N_STO = 20")
# five columns to the right (new I:V) on both row 1 (title/fill band) and
# row 2 (header labels), and opens up D:H for the new N_STO = 2..6 columns.
$ws.Range("D1:H2").EntireColumn.Insert()

# The original D1 formatting (plain, no border/fill) now lives on I1.
# Copy that same "blank" look onto the newly inserted D1:H1 cells so the
# title bar's background/border band extends across all the new columns.
$ws.Range("I1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Fill in the new header labels for row 2 ---------------------------
# Typed in this order so the new shared-string entries are created in the
# same sequence as in the target workbook (N_STO=4,5,6,2,3).
$ws.Range("F2").Value = "N_STO = 4"
$ws.Range("G2").Value = "N_STO = 5"
$ws.Range("H2").Value = "N_STO = 6"
$ws.Range("D2").Value = "N_STO = 2"
$ws.Range("E2").Value = "N_STO = 3"

# --- Row heights (dynamic format) --------------------------------------
$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 15.75

# --- Selection / active cell -------------------------------------------
[void]$ws.Range("E12").Select()
